$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new monthly data point as row 93 (Date / Services / Goods),
# continuing the existing series.
$ws.Range("A93").Value = 45505
$ws.Range("B93").Value = 0.500641946919613
$ws.Range("C93").Value = 0.102185333053162

# The Date column (A2:A92) used a custom "mm/dd/yyyy" number format; switch it
# to Excel's builtin short-date format (numFmtId 14) and make sure the new
# row picks up that same date formatting. Set the format once, then copy it
# (formats only) across the whole column so every cell keeps sharing a single
# underlying style, same as before the edit.
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy()
$ws.Range("A2:A93").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
